# Feria Lagunitas de Puerto Montt - Uva
# Weekly update: insert 3 new price observations (rows 240-242), pushing the
# existing data (old rows 240-322) down to rows 243-325.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 blank rows before the current row 240 (shifts old 240:322 -> 243:325)
$ws.Rows("240:242").Insert()

# New row 240: Red Globe
$ws.Range("A240").Value = 4
$ws.Range("B240").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C240").Value = "Los Lagos"
$ws.Range("D240").Value = 44988
$ws.Range("E240").Value = 10
$ws.Range("F240").Value = "Fruta"
$ws.Range("G240").Value = 100109
$ws.Range("H240").Value = "Uva"
$ws.Range("I240").Value = 100109001
$ws.Range("J240").Value = "Uva"
$ws.Range("K240").Value = "Red Globe"
$ws.Range("L240").Value = "Primera"
$ws.Range("M240").Value = 400
$ws.Range("N240").Value = 15000
$ws.Range("O240").Value = 16000
$ws.Range("P240").Value = 15500
$ws.Range("Q240").Value = "`$/caja 18 kilos"
$ws.Range("R240").Value = "Región de O'Higgins"
$ws.Range("S240").Value = 861
$ws.Range("T240").Value = 18

# New row 241: Rosada pastilla
$ws.Range("A241").Value = 4
$ws.Range("B241").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C241").Value = "Los Lagos"
$ws.Range("D241").Value = 44988
$ws.Range("E241").Value = 10
$ws.Range("F241").Value = "Fruta"
$ws.Range("G241").Value = 100109
$ws.Range("H241").Value = "Uva"
$ws.Range("I241").Value = 100109001
$ws.Range("J241").Value = "Uva"
$ws.Range("K241").Value = "Rosada pastilla"
$ws.Range("L241").Value = "Primera"
$ws.Range("M241").Value = 300
$ws.Range("N241").Value = 18000
$ws.Range("O241").Value = 18000
$ws.Range("P241").Value = 18000
$ws.Range("Q241").Value = "`$/caja 12 kilos"
$ws.Range("R241").Value = "Provincia de Limarí"
$ws.Range("S241").Value = 1500
$ws.Range("T241").Value = 12

# New row 242: Superior Seedless
$ws.Range("A242").Value = 4
$ws.Range("B242").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C242").Value = "Los Lagos"
$ws.Range("D242").Value = 44988
$ws.Range("E242").Value = 10
$ws.Range("F242").Value = "Fruta"
$ws.Range("G242").Value = 100109
$ws.Range("H242").Value = "Uva"
$ws.Range("I242").Value = 100109001
$ws.Range("J242").Value = "Uva"
$ws.Range("K242").Value = "Superior Seedless"
$ws.Range("L242").Value = "Primera"
$ws.Range("M242").Value = 300
$ws.Range("N242").Value = 15000
$ws.Range("O242").Value = 16000
$ws.Range("P242").Value = 15500
$ws.Range("Q242").Value = "`$/caja 18 kilos"
$ws.Range("R242").Value = "Región de O'Higgins"
$ws.Range("S242").Value = 861
$ws.Range("T242").Value = 18
